# Apply "changes to flight selectors with screen scraping" edit.
#
# - Hotel nightly rate for Omni Dallas Hotel / Sonder at Commerce bumped
#   from $110 to $111 (their PRODUCT() totals recalc automatically).
# - The flight-cost column header becomes more specific ("Cost" ->
#   "Cost / ticket").
# - The two scraped flight options (outbound + return) are replaced with
#   a new pair of Qatar Airways/American itineraries (airline, schedule,
#   duration and per-ticket cost all change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode pieces used by the Google-Flights-style "H:MM AM/PM – H:MM AM/PM"
# schedule strings (narrow no-break space before AM/PM, no-break space
# around the en dash).
$nnbsp  = [char]0x202F
$nbsp   = [char]0x00A0
$endash = [char]0x2013

# --- Hotel section -------------------------------------------------------
$ws.Range("B2").Value = 111
$ws.Range("B3").Value = 111

# --- Flight header ---------------------------------------------------------
$ws.Range("E5").Value = "Cost / ticket"

# --- Outbound flight (row 6) ---------------------------------------------
$ws.Range("A6").Value = "Qatar AirwaysAmerican"
$ws.Range("B6").Value = "9:05" + $nnbsp + "AM" + $nbsp + $endash + $nbsp + "10:00" + $nnbsp + "AM+1"
$ws.Range("D6").Value = "36 hr 55 min"
$ws.Range("E6").Value = 1266

# --- Return flight (row 7) -------------------------------------------------
$ws.Range("A7").Value = "Qatar AirwaysAmerican"
$ws.Range("B7").Value = "6:45" + $nnbsp + "PM" + $nbsp + $endash + $nbsp + "9:35" + $nnbsp + "AM+1"
$ws.Range("D7").Value = "26 hr 50 min"
$ws.Range("E7").Value = 1386

# Leave the cursor where the author's session ended up.
$ws.Range("H9").Select()
